$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 1.67
$ws.Range("K2").Value = 4.5
$ws.Range("L2").Value = 1.31
$ws.Range("Q2").Value = 1.7
$ws.Range("AB2").Value = 10.5
$ws.Range("AF2").Value = 11.5

# Row 5
$ws.Range("J5").Value = 2.56
$ws.Range("P5").Value = 1.41
$ws.Range("Q5").Value = 2.7

# Row 6
$ws.Range("P6").Value = 1.93
$ws.Range("Q6").Value = 1.65

# Row 8
$ws.Range("F8").Value = 1.5
$ws.Range("G8").Value = 1.83
$ws.Range("K8").Value = 5.9

# Row 9
$ws.Range("AJ9").Value = 20
$ws.Range("AK9").Value = 23
